# Automatische test-sync: 2025-08-06 20:17:50
# Appends a new mail-log row (row 15) to the "Logs" sheet and the matching
# aggregated row (row 5) to the "Dashboard" sheet, widens the five
# conditional-formatting ranges on "Logs" that tracked the used range, and
# extends the dashboard chart's category/value series so it keeps plotting
# the full data set.

$wb = $excel.ActiveWorkbook

# ---- Logs sheet: append row 15 ----------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(15, 1).Value = "Laat maar weten of er nieuws is"
$logs.Cells.Item(15, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(15, 3).Value = "Testmail #1: Laat maar weten of er nieuws is"
$logs.Cells.Item(15, 4).Value = "Klantenservice / Opvolging"
$logs.Cells.Item(15, 5).Value = "Bedankt, we hebben dit doorgestuurd naar support@bedrijf.nl."
$logs.Cells.Item(15, 6).Value = "2025-08-06 20:17:47"
$logs.Cells.Item(15, 7).Value = "Ja"
$logs.Cells.Item(15, 8).Value = "Ja"
$logs.Cells.Item(15, 9).Value = "Nee"
$logs.Cells.Item(15, 10).Value = "Nee"

# ---- Logs sheet: widen the existing conditional-formatting ranges -----
# (D2:D14, G2:G14, H2:H14, I2:I14, J2:J14 -> same columns, through row 15)
$cfColumns = @("D", "G", "H", "I", "J")
foreach ($col in $cfColumns) {
    $oldRange = $logs.Range($col + "2:" + $col + "14")
    $newRange = $logs.Range($col + "2:" + $col + "15")
    $oldRange.FormatConditions.Item(1).ModifyAppliesToRange($newRange)
}

# ---- Dashboard sheet: append row 5 -------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(5, 1).Value = "Klantenservice / Opvolging"
$dash.Cells.Item(5, 2).Value = 1

# ---- Chart: widen series references from row 4 to row 5 ---------------
$chart = $dash.ChartObjects(1).Chart
$series = $chart.SeriesCollection(1)
$series.XValues = "='Dashboard'!`$A`$2:`$A`$5"
$series.Values = "='Dashboard'!`$B`$2:`$B`$5"
